$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value2 = "Budget Method"
